# 14-Apr-2024: Administrator functions were implemented.
# Adds four new candidate rows (batches 6666 and 5555) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append below the existing table (rows 2-7 already used).
$rows = @(
    @{ Row = 8;  A = 6666; B = "EL54321";  C = "Fung"; D = "Ming Kong"; E = 90156789 },
    @{ Row = 9;  A = 6666; B = "EL666666"; C = "Fung"; D = "Steve";     E = 24484568 },
    @{ Row = 10; A = 5555; B = "EL12345";  C = "Chan"; D = "Tai Man";   E = 98765432 },
    @{ Row = 11; A = 5555; B = "EL98765";  C = "Wong"; D = "Tai Sin";   E = 65432109 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D

    $eCell = $ws.Cells.Item($rowNum, 5)
    $eCell.Value = $r.E
    $eCell.NumberFormat = "@"
}

# Reflect the selection recorded in the saved workbook.
$ws.Range("B10:E11").Select()
